$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.562128560411576
$ws.Range("C2").Value = 2.828655049581922
$ws.Range("D2").Value = 10.56472949907906
$ws.Range("E2").Value = 2.437664035705469
$ws.Range("F2").Value = 2.558139091228547
$ws.Range("G2").Value = 3.657830235493572
$ws.Range("H2").Value = 2.436801249702838
$ws.Range("B3").Value = 2.576937016845974
$ws.Range("C3").Value = 2.883834815612763
$ws.Range("D3").Value = 8.374954051276625
$ws.Range("E3").Value = 2.444871584426221
$ws.Range("F3").Value = 2.571058240846794
$ws.Range("G3").Value = 3.555220259392242
$ws.Range("H3").Value = 2.442838613609139
$ws.Range("B4").Value = 2.526949918108024
$ws.Range("C4").Value = 3.01732006345383
$ws.Range("D4").Value = 9.361953554425824
$ws.Range("E4").Value = 2.41651902216063
$ws.Range("F4").Value = 2.524033221911411
$ws.Range("G4").Value = 3.607583539008169
$ws.Range("H4").Value = 2.416199608207933
$ws.Range("B5").Value = 2.571187619835253
$ws.Range("C5").Value = 3.380719459040112
$ws.Range("D5").Value = 4.940140009134441
$ws.Range("E5").Value = 2.446396964240126
$ws.Range("F5").Value = 2.564825141152768
$ws.Range("G5").Value = 3.534977467695958
$ws.Range("H5").Value = 2.444040844258239
$ws.Range("B6").Value = 2.580933269305367
$ws.Range("C6").Value = 3.364883251118721
$ws.Range("D6").Value = 1.426955281458697
$ws.Range("E6").Value = 2.453820867654344
$ws.Range("F6").Value = 2.574585052914157
$ws.Range("G6").Value = 2.968124494335742
$ws.Range("H6").Value = 2.451383197835503
$ws.Range("B7").Value = 2.555471140536096
$ws.Range("C7").Value = 3.655791957161804
$ws.Range("D7").Value = 1.909465123111566
$ws.Range("E7").Value = 2.445607710036875
$ws.Range("F7").Value = 2.549393159961533
$ws.Range("G7").Value = 3.230341905891171
$ws.Range("H7").Value = 2.44341887967167
$ws.Range("B8").Value = 2.568469745368881
$ws.Range("C8").Value = 1.61754885472335
$ws.Range("D8").Value = 3.162607168728365
$ws.Range("E8").Value = 2.428272206215448
$ws.Range("F8").Value = 2.562226646976647
$ws.Range("G8").Value = 2.048578856136538
$ws.Range("H8").Value = 2.426137845639802
$ws.Range("B9").Value = 2.57466267579235
$ws.Range("C9").Value = 3.03186021144269
$ws.Range("D9").Value = 2.414254612245162
$ws.Range("E9").Value = 2.438785451754285
$ws.Range("F9").Value = 2.569237339612119
$ws.Range("G9").Value = 2.83681288685493
$ws.Range("H9").Value = 2.437293282939875
$ws.Range("B10").Value = 2.253757000171591
$ws.Range("C10").Value = 3.678504460331955
$ws.Range("D10").Value = 4.478426308014347
$ws.Range("E10").Value = 2.350206967072472
$ws.Range("F10").Value = 2.259907216303427
$ws.Range("G10").Value = 3.499139607109537
$ws.Range("H10").Value = 2.354271173204167
$ws.Range("B11").Value = 2.186605345224173
$ws.Range("C11").Value = 3.781325438004655
$ws.Range("D11").Value = 3.624572709325329
$ws.Range("E11").Value = 2.355463786409341
$ws.Range("F11").Value = 2.194203204459443
$ws.Range("G11").Value = 3.48083342614943
$ws.Range("H11").Value = 2.358096505327883
$ws.Range("B12").Value = 1.858851259405834
$ws.Range("C12").Value = 3.573545316688634
$ws.Range("D12").Value = 0.7992521905295669
$ws.Range("E12").Value = 2.232917583795264
$ws.Range("F12").Value = 1.865783552065658
$ws.Range("G12").Value = 3.018377456737963
$ws.Range("H12").Value = 2.233541039399723
$ws.Range("B13").Value = 2.23224328459987
$ws.Range("C13").Value = 3.657686037493022
$ws.Range("D13").Value = 3.085090764804561
$ws.Range("E13").Value = 2.338634367234772
$ws.Range("F13").Value = 2.237924815214063
$ws.Range("G13").Value = 3.337211194759394
$ws.Range("H13").Value = 2.342009859240121
